$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header column (H1) - copy the formatting from the existing
# header cell (G1) so it matches the other header cells (bold, bordered,
# centered), then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New "Save" values for row 2 (default/unformatted numeric cell)
$ws.Range("H2").Value = 0
